$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 68; this shifts the existing rows 68-129 down to 69-130,
# preserving all their original data (dates/prices move down one row each).
$ws.Rows.Item(68).Insert()

# Populate the newly inserted row 68 with the new weekly record.
$ws.Cells.Item(68, 1).Value = 8
$ws.Cells.Item(68, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(68, 3).Value = "Coquimbo"
$ws.Cells.Item(68, 4).Value = 44512
$ws.Cells.Item(68, 5).Value = 4
$ws.Cells.Item(68, 6).Value = 100112037
$ws.Cells.Item(68, 7).Value = "Cebollín"
$ws.Cells.Item(68, 8).Value = "Sin especificar"
$ws.Cells.Item(68, 9).Value = "Primera"
$ws.Cells.Item(68, 10).Value = 3000
$ws.Cells.Item(68, 11).Value = 900
$ws.Cells.Item(68, 12).Value = 1000
$ws.Cells.Item(68, 13).Value = 950
$ws.Cells.Item(68, 14).Value = "`$/paquete 6 unidades"
$ws.Cells.Item(68, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(68, 16).Value = 158
$ws.Cells.Item(68, 17).Value = 6
$ws.Cells.Item(68, 18).Value = "Hortaliza"
